$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 3
$ws.Range("A2").Value = 335208
$ws.Range("B2").Value = 89388
$ws.Range("E2").Value = 1108
$ws.Range("Q2").Value = 503400.1668687108
$ws.Range("R2").Value = 7097087.468780059
$ws.Range("D2").Value = "NT"
$ws.Range("F2").Value = "Harticka"
$ws.Range("G2").Value = "Pelloporus leporinus"
$ws.Range("H2").Value = "(Fr.) Krieglst."
# Row 3 <- old row 4
$ws.Range("A3").Value = 393361
$ws.Range("B3").Value = 89832
$ws.Range("E3").Value = 1209
$ws.Range("Q3").Value = 503312.8822569048
$ws.Range("R3").Value = 7096945.197971026
$ws.Range("D3").Value = "VU"
$ws.Range("F3").Value = "Rynkskinn"
$ws.Range("G3").Value = "Phlebia centrifuga"
$ws.Range("H3").Value = "P.Karst."
# Row 4 <- old row 5
$ws.Range("A4").Value = 1561195
$ws.Range("B4").Value = 89410
$ws.Range("E4").Value = 5432
$ws.Range("Q4").Value = 503459.5351062054
$ws.Range("R4").Value = 7097040.142521238
$ws.Range("D4").Value = "NT"
$ws.Range("F4").Value = "Granticka"
$ws.Range("G4").Value = "Porodaedalea chrysoloma"
$ws.Range("H4").Value = "(Fr.) Fiasson & Niemelä"
# Row 5 <- old row 6
$ws.Range("A5").Value = 883436
$ws.Range("B5").Value = 78570
$ws.Range("E5").Value = 2081
$ws.Range("Q5").Value = 503601.4553612238
$ws.Range("R5").Value = 7097420.731937714
$ws.Range("D5").Value = "NT"
$ws.Range("F5").Value = "Skrovellav"
$ws.Range("G5").Value = "Lobaria scrobiculata"
$ws.Range("H5").Value = "(Scop.) DC."
# Row 6 <- old row 2
$ws.Range("A6").Value = 393360
$ws.Range("B6").Value = 89832
$ws.Range("E6").Value = 1209
$ws.Range("Q6").Value = 503663.2819856483
$ws.Range("R6").Value = 7097147.882865341
$ws.Range("D6").Value = "VU"
$ws.Range("F6").Value = "Rynkskinn"
$ws.Range("G6").Value = "Phlebia centrifuga"
$ws.Range("H6").Value = "P.Karst."
